# Add a new LeetCode tracking entry ("Find the Difference", #389) as the
# next row in the log, matching the formatting of the existing rows above.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 18 has the same "plain" formatting pattern we need for the new row:
# a date cell (style with m/d/yyyy number format) in column A, and plain
# wrap-text cells in B:G with no fill/highlight, and no value in H.
$ws.Range("A18:G18").Copy()
$ws.Range("A25:G25").PasteSpecial(-4122)

# Fill in the new record's values.
$ws.Range("A25").Value = 44098
$ws.Range("B25").Value = 389
$ws.Range("C25").Value = "Find the Difference"
$ws.Range("D25").Value = "Easy"
$ws.Range("E25").Value = "Given two strings s and t which consist of only lowercase letters.`nString t is generated by random shuffling string s and then add one more letter at a random position.`nFind the letter that was added in t."
$ws.Range("F25").Value = 1
$ws.Range("G25").Value = "1.HashMap`n2.Sorting`n3.XOR"

# Match the row height used for this entry.
$ws.Rows.Item(25).RowHeight = 85.5

# Update the view: scroll so row 16 is at the top and select the cell just
# past the newly entered data, as if the user had just finished typing.
$excel.Goto($ws.Range("A16"), $true)
$ws.Range("G26").Select()
